$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 7003.0293
$ws.Range("I98").Value = 6043.727
$ws.Range("K98").Value = 6043.727
$ws.Range("M98").Value = -4545.727
# Row 122
$ws.Range("H122").Value = 7003.0293
$ws.Range("I122").Value = 6043.727
$ws.Range("K122").Value = 18131.181
$ws.Range("M122").Value = -15681.181
# Row 137
$ws.Range("H137").Value = 38880.742
$ws.Range("I137").Value = 1716
$ws.Range("J137").Value = 168957.33
$ws.Range("K137").Value = 5148
$ws.Range("L137").Value = 506871.99
$ws.Range("M137").Value = -2598
$ws.Range("N137").Value = -511971.99

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1962.6154
$ws.Range("I2").Value = 1626.5834
$ws.Range("J2").Value = 5995
$ws.Range("K2").Value = 1626.5834
$ws.Range("L2").Value = 5995
$ws.Range("M2").Value = -1513.5834
$ws.Range("N2").Value = -6221
# Row 45
$ws.Range("H45").Value = 14377.923
$ws.Range("I45").Value = 12342.1
$ws.Range("J45").Value = 21164
$ws.Range("K45").Value = 12342.1
$ws.Range("L45").Value = 21164
$ws.Range("M45").Value = -11965.1
$ws.Range("N45").Value = -21918
# Row 61
$ws.Range("H61").Value = 3127.6191
$ws.Range("I61").Value = 3142.9473
$ws.Range("K61").Value = 3142.9473
$ws.Range("M61").Value = -2930.9473
# Row 110
$ws.Range("H110").Value = 2885.625
$ws.Range("I110").Value = 2869.2856
$ws.Range("K110").Value = 2869.2856
$ws.Range("M110").Value = -824.2856000000002
# Row 116
$ws.Range("H116").Value = 1962.6154
$ws.Range("I116").Value = 1626.5834
$ws.Range("J116").Value = 5995
$ws.Range("K116").Value = 1626.5834
$ws.Range("L116").Value = 5995
$ws.Range("M116").Value = 667.4166
$ws.Range("N116").Value = -10583
# Row 122
$ws.Range("H122").Value = 20409.455
$ws.Range("I122").Value = 27265.5
$ws.Range("J122").Value = 2126.6667
$ws.Range("K122").Value = 81796.5
$ws.Range("L122").Value = 6380.000100000001
$ws.Range("M122").Value = -79346.5
$ws.Range("N122").Value = -11280.0001
# Row 136
$ws.Range("H136").Value = 3127.6191
$ws.Range("I136").Value = 3142.9473
$ws.Range("K136").Value = 9428.841899999999
$ws.Range("M136").Value = -6878.841899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1962.6154
$ws.Range("I3").Value = 1626.5834
$ws.Range("J3").Value = 5995
$ws.Range("K3").Value = 1626.5834
$ws.Range("L3").Value = 5995
$ws.Range("M3").Value = -1512.5834
$ws.Range("N3").Value = -6223
# Row 86
$ws.Range("H86").Value = 1788.6316
$ws.Range("I86").Value = 1502.0769
$ws.Range("K86").Value = 1502.0769
$ws.Range("M86").Value = -379.0769
# Row 88
$ws.Range("H88").Value = 58477
$ws.Range("J88").Value = 58477
$ws.Range("L88").Value = 58477
$ws.Range("N88").Value = -59289
# Row 89
$ws.Range("H89").Value = 1788.6316
$ws.Range("I89").Value = 1502.0769
$ws.Range("K89").Value = 7510.3845
$ws.Range("M89").Value = -1894.3845
# Row 91
$ws.Range("H91").Value = 58477
$ws.Range("J91").Value = 58477
$ws.Range("L91").Value = 58477
$ws.Range("N91").Value = -61285

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6597.75
$ws.Range("I31").Value = 4247.25
$ws.Range("J31").Value = 8948.25
$ws.Range("K31").Value = 4247.25
$ws.Range("L31").Value = 8948.25
$ws.Range("M31").Value = -3952.25
$ws.Range("N31").Value = -9538.25
# Row 34
$ws.Range("H34").Value = 6597.75
$ws.Range("I34").Value = 4247.25
$ws.Range("J34").Value = 8948.25
$ws.Range("K34").Value = 4247.25
$ws.Range("L34").Value = 8948.25
$ws.Range("M34").Value = -4045.25
$ws.Range("N34").Value = -9352.25
# Row 58
$ws.Range("H58").Value = 1443.5
$ws.Range("I58").Value = 1443.5
$ws.Range("K58").Value = 1443.5
$ws.Range("M58").Value = -1240.5
# Row 122
$ws.Range("H122").Value = 5649.875
$ws.Range("I122").Value = 5349.75
$ws.Range("J122").Value = 5950
$ws.Range("K122").Value = 16049.25
$ws.Range("L122").Value = 17850
$ws.Range("M122").Value = -13599.25
$ws.Range("N122").Value = -22750
# Row 136
$ws.Range("H136").Value = 1443.5
$ws.Range("I136").Value = 1443.5
$ws.Range("K136").Value = 4330.5
$ws.Range("M136").Value = -1780.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 61
$ws.Range("H61").Value = 185.07692
$ws.Range("J61").Value = 198.7
$ws.Range("L61").Value = 596.0999999999999
$ws.Range("N61").Value = -1026.1
# Row 92
$ws.Range("H92").Value = 920
$ws.Range("J92").Value = 875
$ws.Range("L92").Value = 2625
$ws.Range("N92").Value = -5121
# Row 113
$ws.Range("H113").Value = 649.125
$ws.Range("I113").Value = 560
$ws.Range("K113").Value = 1680
$ws.Range("M113").Value = 490
# Row 122
$ws.Range("H122").Value = 609.04346
$ws.Range("J122").Value = 672.6842
$ws.Range("L122").Value = 6054.1578
$ws.Range("N122").Value = -10954.1578
# Row 133
$ws.Range("H133").Value = 10635.091
$ws.Range("I133").Value = 5796.4
$ws.Range("K133").Value = 17389.2
$ws.Range("M133").Value = -12329.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = ""
# Row 36
$ws.Range("H36").Value = 6200
$ws.Range("J36").Value = 5600
$ws.Range("L36").Value = 5600
$ws.Range("N36").Value = -6570
# Row 43
$ws.Range("H43").Value = 20202
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 20202
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 20202
$ws.Range("M43").Value = ""
$ws.Range("N43").Value = -20504
# Row 57
$ws.Range("H57").Value = 23583.25
$ws.Range("I57").Value = 9999
$ws.Range("J57").Value = 24818.182
$ws.Range("K57").Value = 9999
$ws.Range("L57").Value = 24818.182
$ws.Range("M57").Value = -9179
$ws.Range("N57").Value = -26458.182
# Row 80
$ws.Range("H80").Value = 10096.071
$ws.Range("I80").Value = 14006.875
$ws.Range("K80").Value = 14006.875
$ws.Range("M80").Value = -13008.875
# Row 83
$ws.Range("H83").Value = 10096.071
$ws.Range("I83").Value = 14006.875
$ws.Range("K83").Value = 70034.375
$ws.Range("M83").Value = -65042.375
# Row 122
$ws.Range("H122").Value = 1116.5
$ws.Range("I122").Value = 1104.9048
$ws.Range("K122").Value = 3314.7144
$ws.Range("M122").Value = -864.7143999999998
# Row 136
$ws.Range("H136").Value = 39395.4
$ws.Range("J136").Value = 39395.4
$ws.Range("L136").Value = 118186.2
$ws.Range("N136").Value = -123286.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 20842290
$ws.Range("I61").Value = 20842290
$ws.Range("K61").Value = 20842290
$ws.Range("M61").Value = -20842088
# Row 68
$ws.Range("H68").Value = 2794.0454
$ws.Range("I68").Value = 2752.6155
$ws.Range("J68").Value = 2853.889
$ws.Range("K68").Value = 2752.6155
$ws.Range("L68").Value = 2853.889
$ws.Range("M68").Value = -2003.6155
$ws.Range("N68").Value = -4351.889
# Row 71
$ws.Range("H71").Value = 2794.0454
$ws.Range("I71").Value = 2752.6155
$ws.Range("J71").Value = 2853.889
$ws.Range("K71").Value = 13763.0775
$ws.Range("L71").Value = 14269.445
$ws.Range("M71").Value = -10019.0775
$ws.Range("N71").Value = -21757.445
# Row 113
$ws.Range("H113").Value = 20842290
$ws.Range("I113").Value = 20842290
$ws.Range("K113").Value = 20842290
$ws.Range("M113").Value = -20840120

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1637
$ws.Range("I122").Value = 1789.909
$ws.Range("K122").Value = 5369.727000000001
$ws.Range("M122").Value = -2919.727000000001
